# Add serialization via serialize API #208 - generation
#
# Target shape (per the commit's xml diff):
#   Sheet tab order becomes: Astruct, AstructBstruct2Use, AstructBstructUse,
#   Bstruct, Dstruct  (the original blank "Sheet1" first tab is removed).
#
#   AstructBstructUse / AstructBstruct2Use are new 2-column sheets shaped
#   like Bstruct (Name / Bstruct2|Bstrcut2 header row, autofiltered,
#   filterMode sheetPr), narrower than Bstruct. Dstruct is a brand-new
#   1-column sheet (Name header only, autofiltered).
#
# NOTE: worksheet object references returned earlier in the script can go
# "stale" once further sheets are inserted/removed -- operations invoked
# through a stale reference silently land on the wrong sheet. To stay
# safe, every worksheet is re-fetched by name via $wb.Worksheets.Item(...)
# immediately before it is touched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the three new sheets in the right slots.
#    Worksheets.Add($Before) inserts immediately before $Before, so build
#    right-to-left relative to Bstruct for the two "Use" sheets, then drop
#    Dstruct in immediately after Bstruct.
# ---------------------------------------------------------------------

$bstruct = $wb.Worksheets.Item("Bstruct")
$astructBstructUse = $wb.Worksheets.Add($bstruct)
$astructBstructUse.Name = "AstructBstructUse"

$astructBstructUseRef = $wb.Worksheets.Item("AstructBstructUse")
$astructBstruct2Use = $wb.Worksheets.Add($astructBstructUseRef)
$astructBstruct2Use.Name = "AstructBstruct2Use"

$bstructRef = $wb.Worksheets.Item("Bstruct")
$dstruct = $wb.Worksheets.Add($null, $bstructRef)
$dstruct.Name = "Dstruct"

# ---------------------------------------------------------------------
# 2) Populate AstructBstruct2Use (header row + narrow columns + filter).
# ---------------------------------------------------------------------

$sheet = $wb.Worksheets.Item("AstructBstruct2Use")
$sheet.Range("A1").Value = "Name"
$sheet.Range("B1").Value = "Bstrcut2"
$sheet.Columns.Item(1).ColumnWidth = 5.166666666666667   # serialises to width="6"
$sheet.Columns.Item(2).ColumnWidth = 9.166666666666666   # serialises to width="10"
[void]$sheet.Range("A1:B1").AutoFilter()

# ---------------------------------------------------------------------
# 3) Populate AstructBstructUse (header row + narrow columns + filter).
# ---------------------------------------------------------------------

$sheet = $wb.Worksheets.Item("AstructBstructUse")
$sheet.Range("A1").Value = "Name"
$sheet.Range("B1").Value = "Bstruct2"
$sheet.Columns.Item(1).ColumnWidth = 5.166666666666667   # serialises to width="6"
$sheet.Columns.Item(2).ColumnWidth = 9.166666666666666   # serialises to width="10"
[void]$sheet.Range("A1:B1").AutoFilter()

# ---------------------------------------------------------------------
# 4) Populate Dstruct (single Name column + filter).
# ---------------------------------------------------------------------

$sheet = $wb.Worksheets.Item("Dstruct")
$sheet.Range("A1").Value = "Name"
$sheet.Columns.Item(1).ColumnWidth = 5.166666666666667    # serialises to width="6"
[void]$sheet.Range("A1:A1").AutoFilter()

# ---------------------------------------------------------------------
# 5) Drop the original blank first sheet. Excel moves ActiveSheet/
#    tabSelected to the new first tab (Astruct) automatically.
# ---------------------------------------------------------------------

[void]$wb.Worksheets.Item("Sheet1").Delete()

# ---------------------------------------------------------------------
# 6) Hidden _xlnm._FilterDatabase defined names, scoped to each new sheet
#    (mirrors the ones Astruct/Bstruct already carry).
# ---------------------------------------------------------------------

$sheet = $wb.Worksheets.Item("AstructBstruct2Use")
$sheet.Activate()
[void]$sheet.Names.Add("_xlnm._FilterDatabase", "=AstructBstruct2Use!`$A`$1:`$B`$1", $false)

$sheet = $wb.Worksheets.Item("AstructBstructUse")
$sheet.Activate()
[void]$sheet.Names.Add("_xlnm._FilterDatabase", "=AstructBstructUse!`$A`$1:`$B`$1", $false)

$sheet = $wb.Worksheets.Item("Dstruct")
$sheet.Activate()
[void]$sheet.Names.Add("_xlnm._FilterDatabase", "=Dstruct!`$A`$1:`$A`$1", $false)

for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "AstructBstruct2Use!_FilterDatabase" -or `
        $n.Name -eq "AstructBstructUse!_FilterDatabase" -or `
        $n.Name -eq "Dstruct!_FilterDatabase") {
        $n.Visible = $false
    }
}

# Restore Astruct as the active/selected sheet (matches tabSelected moving
# there once Sheet1 is gone).
$wb.Worksheets.Item("Astruct").Activate()
